$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.942.47"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").Value = "1.846.37"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  +0.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.64"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4767"
$ws.Range("E7").Value = "  +2.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3665"
$ws.Range("E8").Value = "  +1.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07223"
$ws.Range("E9").Value = "  +1.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9272"
$ws.Range("E10").Value = "  +2.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.69"
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("D13").Value = "1.886.86"
$ws.Range("E13").Value = "  +3.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.316"
$ws.Range("E14").Value = "  +0.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.411"
$ws.Range("E15").Value = "  +1.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.81"
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008637"
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").Value = "26.968.15"
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.55"
$ws.Range("E21").Value = "  +2.41%  "
$ws.Range("E22").Value = "  +0.92%  "
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("E24").Value = "  +0.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.39"
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("E27").Value = "  +1.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.15"
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.957"
$ws.Range("E29").Value = "  +2.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08880"
$ws.Range("E30").Value = "  +0.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.319"
$ws.Range("E31").Value = "  +5.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.171"
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7451"
$ws.Range("E33").Value = "  +0.98%  "
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.723"
$ws.Range("E35").Value = "  -3.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.119"
$ws.Range("E36").Value = "  +3.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01956"
$ws.Range("E37").Value = "  +1.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05269"
$ws.Range("E38").Value = "  +2.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.986"
$ws.Range("E39").Value = "  +2.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5192"
$ws.Range("E40").Value = "  +2.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.991"
$ws.Range("E41").Value = "  +1.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1510"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.200"
$ws.Range("E43").Value = "  +2.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.58"
$ws.Range("E44").Value = "  +5.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4729"
$ws.Range("E45").Value = "  +1.19%  "
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.53"
$ws.Range("E47").Value = "  +3.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.602"
$ws.Range("E48").Value = "  +2.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "65.42"
$ws.Range("E49").Value = "  +2.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06027"
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8865"
$ws.Range("E51").Value = "  +4.26%  "
